$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 463, pushing the existing rows 463-529 down to 464-530.
$ws.Rows.Item(463).EntireRow.Insert()

# Populate the newly inserted row 463 with the new record.
$ws.Cells.Item(463, 1).Value = 10
$ws.Cells.Item(463, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(463, 3).Value = "La Araucanía"
$ws.Cells.Item(463, 4).Value = 44748
$ws.Cells.Item(463, 5).Value = 9
$ws.Cells.Item(463, 6).Value = 100112028
$ws.Cells.Item(463, 7).Value = "Sandia"
$ws.Cells.Item(463, 8).Value = "Sin especificar"
$ws.Cells.Item(463, 9).Value = "Primera"
$ws.Cells.Item(463, 10).Value = 120
$ws.Cells.Item(463, 11).Value = 3125
$ws.Cells.Item(463, 12).Value = 3125
$ws.Cells.Item(463, 13).Value = 3125
$ws.Cells.Item(463, 14).Value = "$/unidad"
$ws.Cells.Item(463, 15).Value = "Brasil"
$ws.Cells.Item(463, 16).Value = 3125
$ws.Cells.Item(463, 17).Value = 1
$ws.Cells.Item(463, 18).Value = "Hortaliza"
